$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)
$ws.Range("D1").EntireColumn.Insert()
$ws.Range("D4").Value = "type in JavaDB"
Write-Output "=== before resize, raw cells ==="
foreach ($colLetter in @("B","C","D","E","F","G","H","I")) {
    $v = $ws.Range($colLetter + "4").Value
    Write-Output ($colLetter + "4 = [" + $v + "]")
}
$lo.Resize($ws.Range("B4:I12"))
Write-Output "=== after resize, ListColumns ==="
foreach ($col in $lo.ListColumns) {
    Write-Output ($col.Index.ToString() + ": " + $col.Name)
}
